$wb = $excel.ActiveWorkbook

# --- plotGrids sheet: remove the leftover/unused "tagPrefix" column ---
$wsGrids = $wb.Worksheets.Item("plotGrids")
$wsGrids.Columns.Item(4).Delete()
$wsGrids.Range("D1").Select() | Out-Null

# --- DataCombined sheet: add xOffsetsUnits / yOffsetsUnits columns ---
$wsData = $wb.Worksheets.Item("DataCombined")

# Insert a new column right after "xOffsets" (col H) for "xOffsetsUnits"
$wsData.Columns.Item(9).Insert()
# Insert a new column right after "yOffsets" (now col J, after the shift above) for "yOffsetsUnits"
$wsData.Columns.Item(11).Insert()

$wsData.Cells.Item(1, 9).Value = "xOffsetsUnits"
$wsData.Cells.Item(1, 11).Value = "yOffsetsUnits"

$wsData.Cells.Item(2, 8).Value = 1
$wsData.Cells.Item(2, 9).Value = "h"

$wsData.Cells.Item(3, 8).Value = 1
$wsData.Cells.Item(3, 9).Value = "min"

$wsData.Activate() | Out-Null
$wsData.Range("I4").Select() | Out-Null
